# Applies the 2021-08-03 (serial 44411) outbreak-path update:
#  - 16 new rows appended to Sheet1's Table1 (rows 31-46)
#  - Date Colours sheet: shift the "Colour Code" palette down one row
#    (new row 6 for 3 Aug, existing colours moved to make room)
#  - "Date Colours" becomes the active/selected tab; Sheet1 scrolls down
#    with H46 as the active cell.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)      # "Sheet1"
$ws2 = $wb.Worksheets.Item(2)      # "Date Colours"

# ---------------------------------------------------------------------
# 1. Append 16 new rows to Table1 on Sheet1 (rows 31 through 46)
# ---------------------------------------------------------------------
$lo = $ws1.ListObjects.Item("Table1")
for ($i = 0; $i -lt 16; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# Keep the Date column's short-date format ("d-mmm") for the new rows.
$ws1.Range("A31:A46").NumberFormat = "d-mmm"
$ws1.Range("A31:A46").Value = 44411

$rows = @(
    @{ Row = 31; B = "T1 f17yo";  C = "T30 m10s"; D = "Taringa"; E = "Indooroopilly State High School";          F = "Indooroopilly State High School" },
    @{ Row = 32; B = "T9 m10s";   C = "T31 m10s"; D = "Taringa"; E = "Brisbane Grammar School";                  F = "Brisbane Grammar School Students" },
    @{ Row = 33; B = "T9 m10s";   C = "T32 m10s"; D = "Taringa"; E = "Brisbane Grammar School";                  F = "Brisbane Grammar School Students" },
    @{ Row = 34; B = "T9 m10s";   C = "T33 m10s"; D = "Taringa"; E = "Brisbane Grammar School";                  F = "Brisbane Grammar School Students" },
    @{ Row = 35; B = "T9 m10s";   C = "T34 m10s"; D = "Taringa"; E = "Brisbane Grammar School";                  F = "Brisbane Grammar School Students" },
    @{ Row = 36; B = "T35 f10s";  C = "T35 f10s"; D = "Taringa"; E = "Brisbane Girls' Grammar School";           F = "Brisbane Girls' Grammar School Students" },
    @{ Row = 37; B = "T35 f10s";  C = "T36";      D = "Taringa"; E = "Brisbane Girls' Grammar School";           F = "Brisbane Girls' Grammar School Teacher" },
    @{ Row = 38; B = "T4 child";  C = "T37 child";D = "Taringa"; E = "Ironside State School";                    F = "Ironside State School" },
    @{ Row = 39; B = "T4 child";  C = "T38 child";D = "Taringa"; E = "Ironside State School";                    F = "Ironside State School" },
    @{ Row = 40; B = "T4 child";  C = "T39 child";D = "Taringa"; E = "Ironside State School";                    F = "Ironside State School" },
    @{ Row = 41; B = "T40";       C = "T40";      D = "Taringa"; E = "Household Contacts";                       F = "Household Contacts related to other cases" },
    @{ Row = 42; B = "T40";       C = "T41";      D = "Taringa"; E = "Household Contacts";                       F = "Household Contacts related to other cases" },
    @{ Row = 43; B = "T40";       C = "T42";      D = "Taringa"; E = "Household Contacts";                       F = "Household Contacts related to other cases" },
    @{ Row = 44; B = "T40";       C = "T43";      D = "Taringa"; E = "Household Contacts";                       F = "Household Contacts related to other cases" },
    @{ Row = 45; B = "T40";       C = "T44";      D = "Taringa"; E = "Household Contacts";                       F = "Household Contacts related to other cases" },
    @{ Row = 46; B = "T40";       C = "T45";      D = "Taringa"; E = "Household Contacts";                       F = "Household Contacts related to other cases" }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws1.Cells.Item($row, 2).Value = $r.B
    $ws1.Cells.Item($row, 3).Value = $r.C
    $ws1.Cells.Item($row, 4).Value = $r.D
    $ws1.Cells.Item($row, 5).Value = $r.E
    $ws1.Cells.Item($row, 6).Value = $r.F
    $ws1.Cells.Item($row, 7).Value = "Delta (B.1.617.2)"
    $ws1.Cells.Item($row, 8).Value = "Wild"
}

# ---------------------------------------------------------------------
# 2. "Date Colours" sheet: new colours for 3 Aug (row 6), shifting the
#    previous "Colour Code" values down one date row.
# ---------------------------------------------------------------------
$ws2.Range("B2").Value = "#f3e8f3"
$ws2.Range("B3").Value = "#e6d1e7"
$ws2.Range("B4").Value = "#dabadb"
$ws2.Range("B5").Value = "#cda4cf"
$ws2.Range("B6").Value = "#c08ec3"

# ---------------------------------------------------------------------
# 3. View state: "Date Colours" becomes the active tab; Sheet1's
#    selection moves to the newly-added last row (H46) and the sheet
#    is scrolled down to keep that area in view.
# ---------------------------------------------------------------------
$ws1.Range("H46").Select()
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1

$ws2.Activate()
$ws2.Range("A1").Select()
